$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (single decimal point) need to be
# forced to remain text, matching the source data which stores them as strings.
$textForceCells = @("D4", "D5", "D6", "D9", "D11", "D12", "D13", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D34", "D36", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cell in $textForceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply all new values (remaining assignments are naturally text already).
$ws.Range("D2").Value = "66.894.07"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "3.475.18"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "599.50"
$ws.Range("E5").Value = "  -3.15%  "
$ws.Range("D6").Value = "147.70"
$ws.Range("E6").Value = "  -4.50%  "
$ws.Range("D7").Value = "3.473.83"
$ws.Range("E7").Value = "  -2.48%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("D11").Value = "7.73"
$ws.Range("E11").Value = "  +4.01%  "
$ws.Range("D12").Value = "0.422"
$ws.Range("E12").Value = "  -3.61%  "
$ws.Range("D13").Value = "0.0000212"
$ws.Range("E13").Value = "  -4.26%  "
$ws.Range("D14").Value = "4.056.46"
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").Value = "31.05"
$ws.Range("E15").Value = "  -6.57%  "
$ws.Range("D16").Value = "3.469.23"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("D17").Value = "66.858.24"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "6.37"
$ws.Range("E19").Value = "  -5.46%  "
$ws.Range("D20").Value = "10.12"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "15.15"
$ws.Range("E21").Value = "  -5.24%  "
$ws.Range("D22").Value = "432.43"
$ws.Range("E22").Value = "  -4.79%  "
$ws.Range("D23").Value = "0.604"
$ws.Range("E23").Value = "  -6.21%  "
$ws.Range("D24").Value = "79.08"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "3.606.90"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  -9.38%  "
$ws.Range("D28").Value = "9.81"
$ws.Range("E28").Value = "  -6.62%  "
$ws.Range("D29").Value = "8.30"
$ws.Range("E29").Value = "  -10.32%  "
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("E31").Value = "  -7.04%  "
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "25.31"
$ws.Range("E34").Value = "  -3.09%  "
$ws.Range("D35").Value = "3.460.67"
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("D36").Value = "5.90"
$ws.Range("E36").Value = "  -7.57%  "
$ws.Range("E37").Value = "  -6.74%  "
$ws.Range("D39").Value = "7.88"
$ws.Range("E39").Value = "  -4.29%  "
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "172.90"
$ws.Range("E41").Value = "  -4.66%  "
$ws.Range("D42").Value = "0.0882"
$ws.Range("E42").Value = "  -3.83%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "2.07"
$ws.Range("E43").Value = "  -12.99%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "5.39"
$ws.Range("E44").Value = "  -4.17%  "
$ws.Range("D45").Value = "0.896"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").Value = "46.33"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "28.80"
$ws.Range("E47").Value = "  -7.16%  "
$ws.Range("D48").Value = "1.24"
$ws.Range("E48").Value = "  -7.42%  "
$ws.Range("D49").Value = "7.44"
$ws.Range("E49").Value = "  -4.58%  "
$ws.Range("D50").Value = "2.38"
$ws.Range("E50").Value = "  -10.01%  "
$ws.Range("D51").Value = "0.966"
$ws.Range("E51").Value = "  -5.00%  "

# Restore default (Normal) style on the forced-text cells so no stray number format remains.
foreach ($cell in $textForceCells) {
    $ws.Range($cell).Style = "Normal"
}
